{"js": "// Add the new \"Collaboration Timeline & Diff Foundation Update\" section to the\n// end of the cumulative requirement-status document. The section reuses the\n// same paragraph layout as every earlier section in the doc: a blank line, a\n// \"---\" divider, a title line, an \"Updated: <date>\" line, a blank line, a\n// tab-separated header row, and a tab-separated data row.\n\nconst FONT_NAME = \"Helvetica Light\";\nconst FONT_SIZE = 12; // half-points 24 == 12pt\n\nfunction escapeXml(s) {\n  return String(s)\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Build a minimal single-package OOXML fragment for one paragraph whose run\n// contains `cells` joined by real <w:tab/> elements (so the result matches\n// Word's own \"type text, press Tab, type text...\" serialization instead of\n// collapsing the tab into a literal \\t character inside one <w:t>).\nfunction tabParagraphOoxml(cells) {\n  const runText = cells\n    .map((c) => `<w:t xml:space=\"preserve\">${escapeXml(c)}</w:t>`)\n    .join(\"<w:tab/>\");\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r><w:rPr>\" +\n    `<w:rFonts w:ascii=\"${FONT_NAME}\" w:hAnsi=\"${FONT_NAME}\" w:cs=\"${FONT_NAME}\"/>` +\n    `<w:sz w:val=\"24\"/><w:sz-cs w:val=\"24\"/>` +\n    \"</w:rPr>\" +\n    runText +\n    \"</w:r></w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Insert a row whose cells are tab-separated into a brand-new paragraph at\n// the very end of the body, preserving <w:tab/> elements on save.\nasync function insertTabRow(body, cells) {\n  const placeholder = body.insertParagraph(\"x\", Word.InsertLocation.end);\n  await context.sync();\n  const rng = placeholder.getRange(Word.RangeLocation.whole);\n  rng.insertOoxml(tabParagraphOoxml(cells), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Insert a plain single-run paragraph (no tabs) at the end of the body.\nfunction insertPlainParagraph(body, text) {\n  const p = body.insertParagraph(text, Word.InsertLocation.end);\n  p.font.name = FONT_NAME;\n  p.font.size = FONT_SIZE;\n  return p;\n}\n\nconst body = context.document.body;\n\ninsertPlainParagraph(body, \"\");\ninsertPlainParagraph(body, \"---\");\ninsertPlainParagraph(\n  body,\n  \"Collaboration Timeline & Diff Foundation Update\"\n);\ninsertPlainParagraph(body, \"Updated: 2026-02-18\");\ninsertPlainParagraph(body, \"\");\nawait context.sync();\n\nawait insertTabRow(body, [\n  \"Module Name\",\n  \"Developed\",\n  \"Partial Developed\",\n  \"Need To Develop\",\n]);\n\nawait insertTabRow(body, [\n  \"Collaboration Timeline / Diff Playback\",\n  \"Added backend version-diff API (added/removed/renamed node deltas), editor collaboration session APIs (heartbeat/list/end), and frontend collaboration panel with active editors + diff summary viewer\",\n  \"Single-store polling model only; no real-time websocket sync or conflict merge\",\n  \"CRDT/OT real-time merge, live cursor presence, conflict resolution UI and permissions workflow\",\n]);\n\nawait context.sync();\n", "ps1": "# Add the new \"Collaboration Timeline & Diff Foundation Update\" section to the\n# end of the cumulative requirement-status document. The section reuses the\n# same paragraph layout as every earlier section in the doc: a blank line, a\n# \"---\" divider, a title line, an \"Updated: <date>\" line, a blank line, a\n# tab-separated header row, and a tab-separated data row.\n\n$d = $word.ActiveDocument\n\nfunction Add-PlainParagraph([string]$text) {\n    $p = $d.Paragraphs.Add()\n    $p.Range.Text = $text\n}\n\nfunction Escape-Xml([string]$s) {\n    $s = $s -replace '&', '&amp;'\n    $s = $s -replace '<', '&lt;'\n    $s = $s -replace '>', '&gt;'\n    return $s\n}\n\n# Insert a row whose cells are tab-separated, using real <w:tab/> elements\n# (rather than a literal tab char inside one run) so the saved OOXML matches\n# Word's own \"type text, press Tab, type text...\" serialization.\n#\n# A fresh \"x\" placeholder paragraph is added first and then its range (minus\n# the trailing paragraph-mark character, which Range always includes) is\n# replaced via InsertXML. Targeting the doc-end range directly is unsafe\n# here: Word collapses/merges the incoming XML into a pre-existing *empty*\n# last paragraph instead of appending after it.\nfunction Add-TabRow([string[]]$cells) {\n    $escaped = @()\n    foreach ($c in $cells) { $escaped += Escape-Xml $c }\n    $runText = [string]::Join(\"<w:tab/>\", ($escaped | ForEach-Object { '<w:t xml:space=\"preserve\">' + $_ + '</w:t>' }))\n\n    $ooxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p><w:r><w:rPr>' +\n        '<w:rFonts w:ascii=\"Helvetica Light\" w:hAnsi=\"Helvetica Light\" w:cs=\"Helvetica Light\"/>' +\n        '<w:sz w:val=\"24\"/><w:sz-cs w:val=\"24\"/>' +\n        '</w:rPr>' + $runText + '</w:r></w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $placeholder = $d.Paragraphs.Add()\n    $placeholder.Range.Text = \"x\"\n    $rng = $d.Range($placeholder.Range.Start, $placeholder.Range.End - 1)\n    [void]$rng.InsertXML($ooxml)\n}\n\nAdd-PlainParagraph \"\"\nAdd-PlainParagraph \"---\"\nAdd-PlainParagraph \"Collaboration Timeline & Diff Foundation Update\"\nAdd-PlainParagraph \"Updated: 2026-02-18\"\nAdd-PlainParagraph \"\"\n\nAdd-TabRow @(\"Module Name\", \"Developed\", \"Partial Developed\", \"Need To Develop\")\n\nAdd-TabRow @(\n    \"Collaboration Timeline / Diff Playback\",\n    \"Added backend version-diff API (added/removed/renamed node deltas), editor collaboration session APIs (heartbeat/list/end), and frontend collaboration panel with active editors + diff summary viewer\",\n    \"Single-store polling model only; no real-time websocket sync or conflict merge\",\n    \"CRDT/OT real-time merge, live cursor presence, conflict resolution UI and permissions workflow\"\n)\n"}
